$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic shuffle of the data in rows 2-6 (columns D, L, M, N, O, P, R, S):
# new row2 <= old row4, new row3 <= old row5, new row4 <= old row6,
# new row5 <= old row2, new row6 <= old row3.
# Row 7 and all other columns (A,B,C,E,F,G,H,I,J,K,Q,T) stay unchanged since they are
# identical across rows 2-6 already.

$ws.Range("D2").Value = 44252
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 13500
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 750

$ws.Range("D3").Value = 44250
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 806

$ws.Range("D4").Value = 44253
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("S4").Value = 806

$ws.Range("D5").Value = 45072
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 16000
$ws.Range("R5").Value = "Provincia de Chacabuco"
$ws.Range("S5").Value = 889

$ws.Range("D6").Value = 45072
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 17000
$ws.Range("P6").Value = 17000
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 944
